# Add the new "VehicleFleet" worksheet as the last sheet (after "Warehouses")
$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "VehicleFleet"

# --- Populate the vehicle fleet data (header row + 27 vehicle rows) ---
$ws.Range("A1").Value = "Warehouse"
$ws.Range("B1").Value = "Plate Nr"
$ws.Range("C1").Value = "Make"
$ws.Range("D1").Value = "Model"
$ws.Range("E1").Value = "Capacity in MT"
$ws.Range("A2").Value = "GITEGA"
$ws.Range("B2").Value = "CD44A95"
$ws.Range("C2").Value = "RENAULT 6X4"
$ws.Range("D2").Value = "350,34"
$ws.Range("E2").Value = 18
$ws.Range("A3").Value = "GITEGA"
$ws.Range("B3").Value = "CD44B02"
$ws.Range("C3").Value = "RENAULT 6X4"
$ws.Range("D3").Value = "350,34"
$ws.Range("E3").Value = 18
$ws.Range("A4").Value = "GITEGA"
$ws.Range("B4").Value = "CD44A89"
$ws.Range("C4").Value = "RENAULT4X4"
$ws.Range("D4").Value = "300,19"
$ws.Range("E4").Value = 8
$ws.Range("A5").Value = "GITEGA"
$ws.Range("B5").Value = "CD44A91"
$ws.Range("C5").Value = "RENAULT 4X4"
$ws.Range("D5").Value = "300,19"
$ws.Range("E5").Value = 8
$ws.Range("A6").Value = "GITEGA"
$ws.Range("B6").Value = "CD44A98"
$ws.Range("C6").Value = "RENAULT 4X4"
$ws.Range("D6").Value = "300,19"
$ws.Range("E6").Value = 8
$ws.Range("A7").Value = "GITEGA"
$ws.Range("B7").Value = "CD44A54"
$ws.Range("C7").Value = "TOYOTA DYNA"
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 3.5
$ws.Range("A8").Value = "GITEGA"
$ws.Range("B8").Value = "E059AIT"
$ws.Range("C8").Value = "TOYOTA PIC-UP"
$ws.Range("D8").Value = "Land cruiser"
$ws.Range("E8").Value = 1.5
$ws.Range("A9").Value = "GITEGA"
$ws.Range("B9").Value = "CD107-98U"
$ws.Range("C9").Value = "TRAILER"
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 15
$ws.Range("A10").Value = "BUJUMBURA"
$ws.Range("B10").Value = "CD44A96"
$ws.Range("C10").Value = "RENAULT 6X4"
$ws.Range("D10").Value = "350,34"
$ws.Range("E10").Value = 18
$ws.Range("A11").Value = "BUJUMBURA"
$ws.Range("B11").Value = "CD44A52"
$ws.Range("C11").Value = "RENAULT 4X4"
$ws.Range("D11").Value = "300,19"
$ws.Range("E11").Value = 8
$ws.Range("A12").Value = "BUJUMBURA"
$ws.Range("B12").Value = "CD44A81"
$ws.Range("C12").Value = "ISUZU"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 4.2
$ws.Range("A13").Value = "BUJUMBURA"
$ws.Range("B13").Value = "CD44A55"
$ws.Range("C13").Value = "ISUZU"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 4.2
$ws.Range("A14").Value = "BUJUMBURA"
$ws.Range("B14").Value = "CD44A86"
$ws.Range("C14").Value = "ISUZU"
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 4.2
$ws.Range("A15").Value = "BUJUMBURA"
$ws.Range("B15").Value = "CD44A87"
$ws.Range("C15").Value = "ISUZU"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 4.2
$ws.Range("A16").Value = "BUJUMBURA"
$ws.Range("B16").Value = "CD44A35"
$ws.Range("C16").Value = "TOYOTA DYNA"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 3.5
$ws.Range("A17").Value = "BUJUMBURA"
$ws.Range("B17").Value = "CD44A25"
$ws.Range("C17").Value = "TOYOTA DYNA"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 3.5
$ws.Range("A18").Value = "BUJUMBURA"
$ws.Range("B18").Value = "CD44A31"
$ws.Range("C18").Value = "TOYOTA  PIC-UP"
$ws.Range("D18").Value = "Land cruiser"
$ws.Range("E18").Value = 1.5
$ws.Range("A19").Value = "BUJUMBURA"
$ws.Range("B19").Value = "E058AIT"
$ws.Range("C19").Value = "TOYOTA  PIC-UP"
$ws.Range("D19").Value = "Land cruiser"
$ws.Range("E19").Value = 1.5
$ws.Range("A20").Value = "NGOZI"
$ws.Range("B20").Value = "CD44A88"
$ws.Range("C20").Value = "RENAULT 6X4"
$ws.Range("D20").Value = "350,34"
$ws.Range("E20").Value = 18
$ws.Range("A21").Value = "NGOZI"
$ws.Range("B21").Value = "CD44A94"
$ws.Range("C21").Value = "RENAULT 6X4"
$ws.Range("D21").Value = "350,34"
$ws.Range("E21").Value = 18
$ws.Range("A22").Value = "NGOZI"
$ws.Range("B22").Value = "CD44B01"
$ws.Range("C22").Value = "RENAULT 6X4"
$ws.Range("D22").Value = "350,34"
$ws.Range("E22").Value = 18
$ws.Range("A23").Value = "NGOZI"
$ws.Range("B23").Value = "CD44A90"
$ws.Range("C23").Value = "RENAULT 4X4"
$ws.Range("D23").Value = "300,19"
$ws.Range("E23").Value = 8
$ws.Range("A24").Value = "NGOZI"
$ws.Range("B24").Value = "CD44A57"
$ws.Range("C24").Value = "RENAULT 4X4"
$ws.Range("D24").Value = "300,19"
$ws.Range("E24").Value = 8
$ws.Range("A25").Value = "NGOZI"
$ws.Range("B25").Value = "CD44A48"
$ws.Range("C25").Value = "TOYOTA DYNA"
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 3.5
$ws.Range("A26").Value = "NGOZI"
$ws.Range("B26").Value = "CD44A43"
$ws.Range("C26").Value = "TOYOTA DYNA"
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 3.5
$ws.Range("A27").Value = "NGOZI"
$ws.Range("B27").Value = "CD44A33"
$ws.Range("C27").Value = "TOYOTA  PIC-UP"
$ws.Range("D27").Value = "Land cruiser"
$ws.Range("E27").Value = 1.5
$ws.Range("A28").Value = "NGOZI"
$ws.Range("B28").Value = "CD107-69U"
$ws.Range("C28").Value = "TRAILER"
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 15

# --- Number formats for data rows (rows 2-28) ---
# Column A (Warehouse) as Text, Column E (Capacity in MT) as a 2-decimal number
$ws.Range("A2:A28").NumberFormat = "@"
$ws.Range("E2:E28").NumberFormat = "0.00"

# --- Header row (row 1) formatting: left-aligned text with thin borders ---
$headerRange = $ws.Range("A1:E1")
$headerRange.HorizontalAlignment = -4131
$headerRange.Borders.Item(9).LineStyle = 1
$headerRange.Borders.Item(11).LineStyle = 1
$headerRange.Columns.Item(1).Borders.Item(7).LineStyle = 0
$headerRange.Columns.Item(5).Borders.Item(10).LineStyle = 0

# --- Turn the range into a formatted Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:E28"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Tabelle35"
$tbl.TableStyle = "TableStyleLight16"

# --- Selection / view state for the new sheet ---
$ws.Range("A1:E28").Select()

Write-Output "VehicleFleet sheet added with $($ws.UsedRange.Rows.Count) rows"
